$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows starting at row 6 (pushes existing rows 6.. down to 9..)
$ws.Range("A6:A8").EntireRow.Insert()

# New data for the inserted rows (0362_Seg1, 0362_Seg2, 0362_Seg3)
$data = @(
    @("0362_Seg1", 4, 4, 0, 1, 9, 5, 38, 23),
    @("0362_Seg2", 3, 0, 0, 0, 2, 3, 21, 24),
    @("0362_Seg3", 7, 5, 0, 1, 12, 3, 61, 14)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 6 + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
